$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price entry for "Alcachofa" (Madrigal / Primera) was recorded.
# Insert a new row at row 84, pushing the existing rows 84-91 down to 85-92,
# then fill the new row 84 with the new entry's data.
$ws.Rows.Item(84).Insert()

$ws.Range("A84").Value = 5
$ws.Range("B84").Value = "Macroferia Regional de Talca"
$ws.Range("C84").Value = "Maule"
$ws.Range("D84").Value = 44783
$ws.Range("D84").NumberFormat = $ws.Range("D85").NumberFormat
$ws.Range("E84").Value = 7
$ws.Range("F84").Value = 100112013
$ws.Range("G84").Value = "Alcachofa"
$ws.Range("H84").Value = "Madrigal"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 300
$ws.Range("K84").Value = 14000
$ws.Range("L84").Value = 14000
$ws.Range("M84").Value = 14000
$ws.Range("N84").Value = '$/caja 40 unidades'
$ws.Range("O84").Value = "Provincia del Elquí"
$ws.Range("P84").Value = 350
$ws.Range("Q84").Value = 40
$ws.Range("R84").Value = "Hortaliza"
